$wb = $excel.ActiveWorkbook

# --- Update the "Gridworld_ndm" sheet's D column values (non-deterministic graph data) ---
$wsNdm = $wb.Worksheets.Item("Gridworld_ndm")

$wsNdm.Range("D2").Value = 12.6
$wsNdm.Range("D3").Value = 13.8
$wsNdm.Range("D4").Value = 9
$wsNdm.Range("D5").Value = 20.100000000000001
$wsNdm.Range("D6").Value = 12.5
$wsNdm.Range("D7").Value = 25.8
$wsNdm.Range("D8").Value = 12
$wsNdm.Range("D9").Value = 13.5
$wsNdm.Range("D10").Value = 6.1
$wsNdm.Range("D11").Value = 10.6
$wsNdm.Range("D12").Value = 1.3
$wsNdm.Range("D13").Value = 4.0999999999999996
$wsNdm.Range("D14").Value = 12.8
$wsNdm.Range("D15").Value = 13.3
$wsNdm.Range("D16").Value = 14.8
$wsNdm.Range("D17").Value = 11.9
$wsNdm.Range("D18").Value = 8.3000000000000007
$wsNdm.Range("D19").Value = 4.5999999999999996
$wsNdm.Range("D20").Value = 5.3
$wsNdm.Range("D21").Value = 12.1
$wsNdm.Range("D22").Value = 5.4
$wsNdm.Range("D23").Value = 12.7
$wsNdm.Range("D24").Value = 6.3

# --- Switch the active/selected sheet from "Cartpole" to "Gridworld_ndm", ---
# --- and move its selection to H12 ---
$wsNdm.Activate()
$wsNdm.Range("H12").Select()
